# Fruta / hortaliza, semanal
# Insert a new weekly price record as row 110 on the "Haba" sheet,
# shifting the existing rows 110-125 down to 111-126.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 110 (pushes old row 110 -> 111, etc.)
$ws.Rows.Item(110).Insert()

# Populate the newly inserted row 110 with the new weekly record
$ws.Cells.Item(110, 1).Value = 5
$ws.Cells.Item(110, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(110, 3).Value = "Maule"
$ws.Cells.Item(110, 4).Value = 45209
$ws.Cells.Item(110, 5).Value = 7
$ws.Cells.Item(110, 6).Value = 100112026
$ws.Cells.Item(110, 7).Value = "Haba"
$ws.Cells.Item(110, 8).Value = "Sin especificar"
$ws.Cells.Item(110, 9).Value = "Primera"
$ws.Cells.Item(110, 10).Value = 500
$ws.Cells.Item(110, 11).Value = 9000
$ws.Cells.Item(110, 12).Value = 9000
$ws.Cells.Item(110, 13).Value = 9000
$ws.Cells.Item(110, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(110, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(110, 16).Value = 360
$ws.Cells.Item(110, 17).Value = 25
$ws.Cells.Item(110, 18).Value = "Hortaliza"
